$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This text is shown on the Overview sheet (zh-cn/de-de status columns)
#    and on each language sheet's Status column, for both data rows.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Report the handback: fill in "Latest Target File" (F) and
#    "Latest Handback File" (G) columns, and set the real
#    "Latest Handback DateTime" (H), replacing the 0001-01-01 placeholder.
# ---------------------------------------------------------------------------

$mdFile    = "447830c7-cebb-489f-ba28-ec609559dc98.md"
$zhCnXlf   = "447830c7-cebb-489f-ba28-ec609559dc98.f28090322b799169d72c656f0eb2c2578504b8d6.zh-cn.xlf"
$deDeXlf   = "447830c7-cebb-489f-ba28-ec609559dc98.f28090322b799169d72c656f0eb2c2578504b8d6.de-de.xlf"

$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/92d413a2f2ae5de83552fd4e6703d44efd3decad/e2e/447830c7-cebb-489f-ba28-ec609559dc98.md"
$zhCnXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/27e90114699d5a54e0f64695e1a7dbbb1ca7c292/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/447830c7-cebb-489f-ba28-ec609559dc98.f28090322b799169d72c656f0eb2c2578504b8d6.zh-cn.xlf"
$deDeXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/562c7cd4fe873b4a0f8d8f3ec23f58883c5e91da/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/447830c7-cebb-489f-ba28-ec609559dc98.f28090322b799169d72c656f0eb2c2578504b8d6.de-de.xlf"

$mdDisplayRow3 = "ffff865736d8-fcea-42ea-a321-97636be0831d.md"

# -- zh-cn sheet -------------------------------------------------------------

$wsZhCn.Range("H2").Value = "2016-03-23 04:43:56"
$wsZhCn.Range("H3").Value = "2016-03-23 04:43:56"

# Rebuild the hyperlinks in row/column order so relationship ids line up the
# way Excel lays them out: A2, D2, F2, G2, A3, D3, F3, G3
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, "", "", $mdFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhCnXlfUrl, "", "", $zhCnXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $mdUrl, "", "", $mdFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhCnXlfUrl, "", "", $zhCnXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl, "", "", $mdDisplayRow3)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhCnXlfUrl, "", "", $zhCnXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $mdUrl, "", "", $mdFile)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhCnXlfUrl, "", "", $zhCnXlf)

$wsZhCn.Range("F2").Font.Underline = 2
$wsZhCn.Range("F2").Font.Color = 6603053
$wsZhCn.Range("G2").Font.Underline = 2
$wsZhCn.Range("G2").Font.Color = 6603053
$wsZhCn.Range("F3").Font.Underline = 2
$wsZhCn.Range("F3").Font.Color = 6603053
$wsZhCn.Range("G3").Font.Underline = 2
$wsZhCn.Range("G3").Font.Color = 6603053

# -- de-de sheet ---------------------------------------------------------

$wsDeDe.Range("H2").Value = "2016-03-23 04:44:13"
$wsDeDe.Range("H3").Value = "2016-03-23 04:44:13"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, "", "", $mdFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deDeXlfUrl, "", "", $deDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $mdUrl, "", "", $mdFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deDeXlfUrl, "", "", $deDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl, "", "", $mdDisplayRow3)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deDeXlfUrl, "", "", $deDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $mdUrl, "", "", $mdFile)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deDeXlfUrl, "", "", $deDeXlf)

$wsDeDe.Range("F2").Font.Underline = 2
$wsDeDe.Range("F2").Font.Color = 6603053
$wsDeDe.Range("G2").Font.Underline = 2
$wsDeDe.Range("G2").Font.Color = 6603053
$wsDeDe.Range("F3").Font.Underline = 2
$wsDeDe.Range("F3").Font.Color = 6603053
$wsDeDe.Range("G3").Font.Underline = 2
$wsDeDe.Range("G3").Font.Color = 6603053

$wb.Save()
